$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("solve time") updates reflecting changed MP time limit
$ws.Range("D2").Value = 3627.08570026
$ws.Range("D3").Value = 3643.359197051
$ws.Range("D4").Value = 3625.166957107
$ws.Range("D5").Value = 3626.97548697
$ws.Range("D6").Value = 3648.966281197
$ws.Range("D7").Value = 3627.445946476
$ws.Range("D8").Value = 3630.70043795
$ws.Range("D9").Value = 3628.673196055
$ws.Range("D10").Value = 3629.689509792
$ws.Range("D11").Value = 3632.019231597

# Corrected error in fixed recourse data (objective/gap) for row 7 and 10
$ws.Range("B7").Value = -571.8223658731081
$ws.Range("C7").Value = 0.19879039647430619
$ws.Range("C10").Value = 0.7180756032894822
